# Update "gh-pages" generated output workbook (杭州-漫展信息) to the data
# scraped at commit 456a3b4.
#
# Sheet "展览" (sheet1): most rows get a refreshed "want to go" count (col F);
#   a brand-new event is inserted as row 42 (everything from row 42 on
#   shifts down by one), and the old rows 42-44 (now 43-45) also pick up a
#   couple of updated F values.
# Sheet "演出" (sheet2): a single F-column refresh.
# Sheet "全部类型" (sheet4): the same kind of F-column refreshes as sheet1,
#   but this sheet is NOT restructured - no row is inserted here.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$f1 = @{
    3 = 248
    4 = 251
    5 = 1790
    7 = 326
    8 = 514
    9 = 4719
    10 = 56
    14 = 1288
    17 = 2988
    18 = 1818
    19 = 111
    22 = 46
    24 = 935
    25 = 307
    27 = 2710
    28 = 1021
    29 = 2464
    30 = 250
    31 = 1334
    32 = 716
    34 = 890
    35 = 432
    36 = 1133
    37 = 927
    38 = 1189
    39 = 21
    40 = 868
    41 = 530
}
foreach ($row in $f1.Keys) {
    $ws1.Cells.Item($row, 6).Value = $f1[$row]
}

# Insert the new event as row 42, pushing the old rows 42-44 down to 43-45.
$ws1.Rows.Item(42).Insert()

# Carry the formatting of column A (bold / bordered / centered) onto the
# newly-inserted row's A cell, matching the rest of the table.
$ws1.Range("A41").Copy()
$ws1.Range("A42").PasteSpecial(-4122)   # xlPasteFormats
$ws1.Application.CutCopyMode = $false

$ws1.Cells.Item(42, 1).Value = 41
$ws1.Cells.Item(42, 2).Value = "'2024-08-04"
$ws1.Cells.Item(42, 3).Value = "杭州·梦漫星河动漫嘉年华·赵路专场"
$ws1.Cells.Item(42, 4).Value = "阳城路雅澳杭州电商产业园西侧约200米 杭州大会展中心"
$ws1.Cells.Item(42, 5).Value = "2024.08.04 11:40-08.04 17:00"
$ws1.Cells.Item(42, 6).Value = 64
$ws1.Cells.Item(42, 7).Value = 238
$ws1.Cells.Item(42, 8).Value = "https://show.bilibili.com/platform/detail.html?id=86221"
$ws1.Cells.Item(42, 9).Value = "//i1.hdslb.com/bfs/openplatform/202405/2padflbr1716372780297.jpeg"

# Resequence column A for the rows that shifted down, and refresh the two
# F values that also changed for those rows.
$ws1.Cells.Item(43, 1).Value = 42
$ws1.Cells.Item(44, 1).Value = 43
$ws1.Cells.Item(45, 1).Value = 44

$ws1.Cells.Item(43, 6).Value = 364
$ws1.Cells.Item(45, 6).Value = 3492

# ---------------------------------------------------------------------
# Sheet "演出"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(11, 6).Value = 886

# ---------------------------------------------------------------------
# Sheet "全部类型"
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$f4 = @{
    3 = 248
    4 = 251
    6 = 1790
    8 = 326
    9 = 514
    10 = 4719
    11 = 56
    14 = 1288
    15 = 2988
    17 = 1818
    18 = 111
    21 = 886
    24 = 46
    26 = 935
    27 = 307
    28 = 2710
    31 = 1021
    32 = 2464
    33 = 1334
    34 = 716
    37 = 890
    38 = 1133
    39 = 927
    41 = 1189
    42 = 868
    43 = 530
    44 = 364
    49 = 3492
}
foreach ($row in $f4.Keys) {
    $ws4.Cells.Item($row, 6).Value = $f4[$row]
}
